# Update the SQL queries on Sheet1 so that the LEFT JOIN clauses use the
# renamed key columns (study_id / participant_id) instead of the old
# generic "id" columns, matching the updated C3DC data model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The block of LEFT JOIN clauses shared verbatim by every query cell.
$oldJoin = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`n" + `
    "LEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`n" + `
    "LEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`n" + `
    "LEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`n" + `
    "LEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`n" + `
    "LEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newJoin = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`n" + `
    "LEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`n" + `
    "LEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`n" + `
    "LEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`n" + `
    "LEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`n" + `
    "LEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

# Cells C2, B2, B3, B4, B5, B6, B7 each hold one full SQL query that
# contains exactly one copy of the JOIN block above.
$cellAddresses = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellAddresses) {
    $rng = $ws.Range($addr)
    $text = $rng.Value2
    if ($text -ne $null -and $text.Contains($oldJoin)) {
        $rng.Value2 = $text.Replace($oldJoin, $newJoin)
    }
}
